$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'26.909.79"
$ws.Range('E2').Value = '  -1.89%  '

# Row 3
$ws.Range('D3').Value = "'1.831.44"
$ws.Range('E3').Value = '  -2.14%  '

# Row 4
$ws.Range('D4').Value = "'1.009"
$ws.Range('E4').Value = '  +0.61%  '

# Row 5
$ws.Range('D5').Value = "'311.82"
$ws.Range('E5').Value = '  -1.22%  '

# Row 6
$ws.Range('D6').Value = "'1.008"
$ws.Range('E6').Value = '  +0.60%  '

# Row 7
$ws.Range('D7').Value = "'0.4604"
$ws.Range('E7').Value = '  -1.44%  '

# Row 8
$ws.Range('D8').Value = "'0.3659"
$ws.Range('E8').Value = '  -2.14%  '

# Row 9
$ws.Range('D9').Value = "'0.07208"
$ws.Range('E9').Value = '  -2.62%  '

# Row 10
$ws.Range('D10').Value = "'0.8788"
$ws.Range('E10').Value = '  -1.55%  '

# Row 11
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = "'0.07851"
$ws.Range('E11').Value = '  -1.31%  '

# Row 12
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').Value = "'19.65"
$ws.Range('E12').Value = '  -2.66%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = "'1.899.27"
$ws.Range('E13').Value = '  +1.73%  '

# Row 14
$ws.Range('D14').Value = "'5.337"
$ws.Range('E14').Value = '  -2.14%  '

# Row 15
$ws.Range('D15').Value = "'6.433"
$ws.Range('E15').Value = '  -3.05%  '

# Row 16
$ws.Range('D16').Value = "'89.92"
$ws.Range('E16').Value = '  -3.32%  '

# Row 17
$ws.Range('D17').Value = "'1.009"
$ws.Range('E17').Value = '  +0.63%  '

# Row 18
$ws.Range('D18').Value = "'0.000008738"
$ws.Range('E18').Value = '  -2.70%  '

# Row 19
$ws.Range('E19').Value = '  +0.37%  '

# Row 20
$ws.Range('D20').Value = "'26.955.57"
$ws.Range('E20').Value = '  -1.82%  '

# Row 21
$ws.Range('D21').Value = "'14.51"
$ws.Range('E21').Value = '  -3.14%  '

# Row 22
$ws.Range('D22').Value = "'5.005"
$ws.Range('E22').Value = '  -3.45%  '

# Row 23
$ws.Range('D23').Value = "'10.44"
$ws.Range('E23').Value = '  -1.89%  '

# Row 24
$ws.Range('D24').Value = "'2.004"
$ws.Range('E24').Value = '  +6.64%  '

# Row 25
$ws.Range('D25').Value = "'150.47"
$ws.Range('E25').Value = '  -1.50%  '

# Row 26
$ws.Range('D26').Value = "'18.27"
$ws.Range('E26').Value = '  -1.91%  '

# Row 27
$ws.Range('D27').Value = "'1.998"
$ws.Range('E27').Value = '  -5.17%  '

# Row 28
$ws.Range('D28').Value = "'114.51"
$ws.Range('E28').Value = '  -2.66%  '

# Row 29
$ws.Range('D29').Value = "'4.933"
$ws.Range('E29').Value = '  -4.96%  '

# Row 30
$ws.Range('D30').Value = "'0.08825"
$ws.Range('E30').Value = '  -1.01%  '

# Row 31
$ws.Range('D31').Value = "'3.122"
$ws.Range('E31').Value = '  +5.16%  '

# Row 32
$ws.Range('D32').Value = "'0.7600"
$ws.Range('E32').Value = '  +0.29%  '

# Row 33
$ws.Range('D33').Value = "'4.456"
$ws.Range('E33').Value = '  -1.38%  '

# Row 34
$ws.Range('D34').Value = "'1.135"
$ws.Range('E34').Value = '  -2.69%  '

# Row 35
$ws.Range('D35').Value = "'2.662"
$ws.Range('E35').Value = '  -0.25%  '

# Row 36
$ws.Range('D36').Value = "'1.092"
$ws.Range('E36').Value = '  +0.51%  '

# Row 37
$ws.Range('E37').Value = '  -1.72%  '

# Row 38
$ws.Range('D38').Value = "'0.05155"
$ws.Range('E38').Value = '  -2.81%  '

# Row 39
$ws.Range('D39').Value = "'2.927"
$ws.Range('E39').Value = '  -2.33%  '

# Row 40
$ws.Range('D40').Value = "'6.943"
$ws.Range('E40').Value = '  -3.77%  '

# Row 41
$ws.Range('D41').Value = "'0.4994"
$ws.Range('E41').Value = '  -4.88%  '

# Row 42
$ws.Range('D42').Value = "'0.1598"
$ws.Range('E42').Value = '  -3.19%  '

# Row 43
$ws.Range('D43').Value = "'8.326"
$ws.Range('E43').Value = '  -0.45%  '

# Row 44
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = "'1.008"
$ws.Range('E44').Value = '  +0.62%  '

# Row 45
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = "'0.4668"
$ws.Range('E45').Value = '  -5.42%  '

# Row 46
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = "'10.22"
$ws.Range('E46').Value = '  -1.84%  '

# Row 47
$ws.Range('D47').Value = "'102.55"
$ws.Range('E47').Value = '  -1.33%  '

# Row 48
$ws.Range('D48').Value = "'1.608"
$ws.Range('E48').Value = '  -3.06%  '

# Row 49
$ws.Range('D49').Value = "'0.06125"
$ws.Range('E49').Value = '  -2.41%  '

# Row 50
$ws.Range('D50').Value = "'64.97"
$ws.Range('E50').Value = '  -1.77%  '

# Row 51
$ws.Range('D51').Value = "'36.22"
$ws.Range('E51').Value = '  -2.72%  '
